$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2016Stock")

# Remove the leftover single-cell title row at the top of the sheet; this
# shifts every subsequent row up by one (headers -> row 1, data follows,
# trailing blank/"total" rows shift accordingly).
$ws.Rows.Item(1).Delete()

# Make "2016Stock" the active sheet/tab (previously "2020Demo" was active).
$ws.Activate()

# Select the header row, matching the selection state left behind in the
# saved file.
$ws.Rows.Item(1).Select()
